$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 6 (SanityTest): update username, password, and user detail
$ws.Range("B6").Value = "test123@abv.bg"

# Force the password cell to stay text (it looks numeric) without leaving a
# lingering custom number-format/style on the cell.
$ws.Range("C6").NumberFormat = "@"
$ws.Range("C6").Value = "123456"
$ws.Range("C6").Style = "Normal"

$ws.Range("F6").Value = "Test Testing"

# Row 7 (RegistrationDDT): update username
$ws.Range("B7").Value = "testnew1@abv.bg"
